$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cache enum render value: uppercase the rendered enum display strings
$ws.Range("A2").Value = "BUILDYOURDREAM"
$ws.Range("B2").Value = "ANTIQUEWHITE"

# Widen columns A and B to fit the new (longer/uppercased) cached values.
# (ColumnWidth is specified in characters of the Normal-style font; the
# underlying engine stores/quantizes width in 1/6-character pixel units,
# so these values are chosen to land as close as possible to the target
# stored widths of 19 and 15 characters respectively.)
$ws.Columns.Item(1).ColumnWidth = 18.8333333333333
$ws.Columns.Item(2).ColumnWidth = 14.8333333333333
